$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# Insert a new row above row 12 ("Get flag" shifts down along with everything
# below it). Excel's Rows.Insert shifts formulas/refs automatically.
$ws.Rows("12:12").Insert()

# --- New row 12: "Checkpoint 2626" ---
$ws.Range("A12:D12").Style = "Normal"
$ws.Range("A12").Value = "Checkpoint 2626"
$ws.Range("B12").Value = 3718
$ws.Range("C12").Value = 4192
$ws.Range("D12").Formula = "=IF(B12 >  0,C12-B12, 0)"

# --- Row 13 ("Get flag", previously row 12): new timer values ---
$ws.Range("B13").Value = 4153
$ws.Range("C13").Value = 4627

# --- Row 14 ("Black screen", previously row 13): new timer values ---
$ws.Range("B14").Value = 4671
$ws.Range("C14").Value = 5145

# --- Row 15 ("map - first move", previously row 14): new C value only ---
$ws.Range("C15").Value = 5738

# --- Row 17 ("enter 1-2", previously row 16): drop B, update C ---
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = 5937

# --- Row 18 ("End 1-2", previously row 17): drop both B and C ---
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()

# Restore the copy-down formatting style (col A/D use style 16) that Insert
# carried from row 11, matching the rest of the table rows.
$ws.Range("A12:D12").Style = "Normal"

$ws.Range("B15").Select()
